$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before current row 9, pushing existing rows 9-10 down to 11-12.
$ws.Rows("9:10").Insert()

# New row 9: updated entry previously at old row 9 (now moved to row 11), with new values
$ws.Cells.Item(9, 1).Value = 7
$ws.Cells.Item(9, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(9, 3).Value = "Ñuble"
$ws.Cells.Item(9, 4).Value = 45142
$ws.Cells.Item(9, 5).Value = 16
$ws.Cells.Item(9, 6).Value = "Fruta"
$ws.Cells.Item(9, 7).Value = 100102
$ws.Cells.Item(9, 8).Value = "Cítricos"
$ws.Cells.Item(9, 9).Value = 100102006
$ws.Cells.Item(9, 10).Value = "Pomelo"
$ws.Cells.Item(9, 11).Value = "Start Ruby"
$ws.Cells.Item(9, 12).Value = "Primera"
$ws.Cells.Item(9, 13).Value = 30
$ws.Cells.Item(9, 14).Value = 15000
$ws.Cells.Item(9, 15).Value = 15000
$ws.Cells.Item(9, 16).Value = 15000
$ws.Cells.Item(9, 17).Value = "`$/caja 14 kilos empedrada"
$ws.Cells.Item(9, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(9, 19).Value = 1071
$ws.Cells.Item(9, 20).Value = 14

# New row 10: another new entry
$ws.Cells.Item(10, 1).Value = 7
$ws.Cells.Item(10, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(10, 3).Value = "Ñuble"
$ws.Cells.Item(10, 4).Value = 45142
$ws.Cells.Item(10, 5).Value = 16
$ws.Cells.Item(10, 6).Value = "Fruta"
$ws.Cells.Item(10, 7).Value = 100102
$ws.Cells.Item(10, 8).Value = "Cítricos"
$ws.Cells.Item(10, 9).Value = 100102006
$ws.Cells.Item(10, 10).Value = "Pomelo"
$ws.Cells.Item(10, 11).Value = "Start Ruby"
$ws.Cells.Item(10, 12).Value = "Primera"
$ws.Cells.Item(10, 13).Value = 30
$ws.Cells.Item(10, 14).Value = 14000
$ws.Cells.Item(10, 15).Value = 14000
$ws.Cells.Item(10, 16).Value = 14000
$ws.Cells.Item(10, 17).Value = "`$/caja 14 kilos granel"
$ws.Cells.Item(10, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(10, 19).Value = 1000
$ws.Cells.Item(10, 20).Value = 14

# Apply the date style (style index 2 in original, equiv to number format) to D9 and D10
$ws.Cells.Item(9, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(10, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
